# Insert a new data row at row 147, pushing the existing rows 147-248
# down to 148-249 (dimension grows from A1:T248 to A1:T249), then
# populate the newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(147).Insert()

$ws.Range("A147").Value = 10
$ws.Range("B147").Value = "Vega Modelo de Temuco"
$ws.Range("C147").Value = "La Araucanía"
$ws.Range("D147").Value = 44904
$ws.Range("E147").Value = 9
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100103
$ws.Range("H147").Value = "Frutos de hueso (carozo)"
$ws.Range("I147").Value = 100103002
$ws.Range("J147").Value = "Ciruela"
$ws.Range("K147").Value = "Angeleno"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 65
$ws.Range("N147").Value = 28000
$ws.Range("O147").Value = 28000
$ws.Range("P147").Value = 28000
$ws.Range("Q147").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R147").Value = "Región de O'Higgins"
$ws.Range("S147").Value = 1556
$ws.Range("T147").Value = 18
